# Generate Report for Handback
# -----------------------------------------------------------------------
# This script fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for the two localized-language sheets
# (zh-cn and de-de) now that both languages have been handed back in sync
# with en-US, adds hyperlinks on the newly filled "Latest Target File"
# cells, flips the shared "Status" text from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is used, and widens a
# few columns that now need to show the longer text / file names.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa8987de776f04a1ce75eadb0e196dd698ffdd28"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status text + wider zh-cn / de-de columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-25 20:35:23"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-25 20:35:23"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(10).ColumnWidth = 39.2

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "$repoBase/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "$repoBase/zh-cn/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "$repoBase/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "$repoBase/zh-cn/e2e/a.md", "", "", "a.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-25 20:35:30"

$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-25 20:35:30"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(10).ColumnWidth = 39.2

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "$repoBase/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "$repoBase/de-de/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "$repoBase/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "$repoBase/de-de/e2e/a.md", "", "", "a.md")
